$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2742.4443
$ws.Range("J19").Value = 597.4286
$ws.Range("L19").Value = 597.4286
$ws.Range("N19").Value = -947.4286
$ws.Range("H42").Value = 147.5
$ws.Range("J42").Value = 163.33333
$ws.Range("L42").Value = 489.99999
$ws.Range("N42").Value = -949.99999
$ws.Range("H58").Value = 2377.1428
$ws.Range("I58").Value = 546.6667
$ws.Range("J58").Value = 3750
$ws.Range("K58").Value = 1640.0001
$ws.Range("L58").Value = 11250
$ws.Range("M58").Value = -1490.0001
$ws.Range("N58").Value = -11550
$ws.Range("H64").Value = 4000
$ws.Range("J64").Value = 4000
$ws.Range("L64").Value = 4000
$ws.Range("N64").Value = -4496
$ws.Range("H67").Value = 4000
$ws.Range("J67").Value = 4000
$ws.Range("L67").Value = 4000
$ws.Range("N67").Value = -5716
$ws.Range("H97").Value = 1218
$ws.Range("J97").Value = 1218
$ws.Range("L97").Value = 3654
$ws.Range("N97").Value = -4646
$ws.Range("H116").Value = 3813.7058
$ws.Range("I116").Value = 2091.889
$ws.Range("K116").Value = 2091.889
$ws.Range("M116").Value = 1350.111
$ws.Range("H129").Value = 1074.6428
$ws.Range("I129").Value = 439.66666
$ws.Range("J129").Value = 1196.234
$ws.Range("K129").Value = 1318.99998
$ws.Range("L129").Value = 3588.702
$ws.Range("M129").Value = 3681.00002
$ws.Range("N129").Value = -13588.702
$ws.Range("H132").Value = 1967.9375
$ws.Range("I132").Value = 2069.568
$ws.Range("J132").Value = 850
$ws.Range("K132").Value = 6208.704000000001
$ws.Range("L132").Value = 2550
$ws.Range("M132").Value = -3678.704000000001
$ws.Range("N132").Value = -7610
$ws.Range("H135").Value = 21746432
$ws.Range("I135").Value = 719.1111
$ws.Range("K135").Value = 6471.9999
$ws.Range("M135").Value = -3936.9999
$ws.Range("H137").Value = 1513.7
$ws.Range("J137").Value = 2116.6667
$ws.Range("L137").Value = 6350.000100000001
$ws.Range("N137").Value = -11450.0001
$ws.Range("H138").Value = 2101.5334
$ws.Range("I138").Value = 1635.5454
$ws.Range("J138").Value = 2294.9622
$ws.Range("K138").Value = 4906.6362
$ws.Range("L138").Value = 6884.8866
$ws.Range("M138").Value = 233.3638000000001
$ws.Range("N138").Value = -17164.8866
$ws.Range("H141").Value = 1427.1351
$ws.Range("I141").Value = 1016.4516
$ws.Range("J141").Value = 3549
$ws.Range("K141").Value = 3049.3548
$ws.Range("L141").Value = 10647
$ws.Range("M141").Value = 2130.6452
$ws.Range("N141").Value = -21007

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6068.061
$ws.Range("I32").Value = 4575.8184
$ws.Range("K32").Value = 4575.8184
$ws.Range("M32").Value = -4288.8184
$ws.Range("H74").Value = 50000880
$ws.Range("I74").Value = 83333930
$ws.Range("J74").Value = 1312.375
$ws.Range("K74").Value = 83333930
$ws.Range("L74").Value = 1312.375
$ws.Range("M74").Value = -83333056
$ws.Range("N74").Value = -3060.375
$ws.Range("H77").Value = 50000880
$ws.Range("I77").Value = 83333930
$ws.Range("J77").Value = 1312.375
$ws.Range("K77").Value = 416669650
$ws.Range("L77").Value = 6561.875
$ws.Range("M77").Value = -416665282
$ws.Range("N77").Value = -15297.875
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H122").Value = 2167.9688
$ws.Range("I122").Value = 1802.7858
$ws.Range("J122").Value = 4724.25
$ws.Range("K122").Value = 5408.357400000001
$ws.Range("L122").Value = 14172.75
$ws.Range("M122").Value = -2958.357400000001
$ws.Range("N122").Value = -19072.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H128").Value = 2999.8
$ws.Range("I128").Value = 2999.8
$ws.Range("K128").Value = 8999.400000000001
$ws.Range("M128").Value = -6509.400000000001
$ws.Range("H134").Value = 4791.269
$ws.Range("I134").Value = 5336.0454
$ws.Range("K134").Value = 16008.1362
$ws.Range("M134").Value = -13473.1362

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14586.725
$ws.Range("I31").Value = 36765.89
$ws.Range("J31").Value = 4606.1
$ws.Range("K31").Value = 36765.89
$ws.Range("L31").Value = 4606.1
$ws.Range("M31").Value = -36470.89
$ws.Range("N31").Value = -5196.1
$ws.Range("H34").Value = 14586.725
$ws.Range("I34").Value = 36765.89
$ws.Range("J34").Value = 4606.1
$ws.Range("K34").Value = 36765.89
$ws.Range("L34").Value = 4606.1
$ws.Range("M34").Value = -36563.89
$ws.Range("N34").Value = -5010.1
$ws.Range("H58").Value = 9829.200000000001
$ws.Range("I58").Value = 691.6053000000001
$ws.Range("J58").Value = 30254.412
$ws.Range("K58").Value = 691.6053000000001
$ws.Range("L58").Value = 30254.412
$ws.Range("M58").Value = -488.6053000000001
$ws.Range("N58").Value = -30660.412
$ws.Range("H134").Value = 808.7368
$ws.Range("I134").Value = 697.17645
$ws.Range("J134").Value = 1757
$ws.Range("K134").Value = 2091.52935
$ws.Range("L134").Value = 5271
$ws.Range("M134").Value = 443.4706499999998
$ws.Range("N134").Value = -10341
$ws.Range("H136").Value = 9829.200000000001
$ws.Range("I136").Value = 691.6053000000001
$ws.Range("J136").Value = 30254.412
$ws.Range("K136").Value = 2074.8159
$ws.Range("L136").Value = 90763.236
$ws.Range("M136").Value = 475.1840999999999
$ws.Range("N136").Value = -95863.236

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 200
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 600
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -824
$ws.Range("H97").Value = 3079.6
$ws.Range("I97").Value = 200
$ws.Range("J97").Value = 3799.5
$ws.Range("K97").Value = 600
$ws.Range("L97").Value = 11398.5
$ws.Range("M97").Value = -104
$ws.Range("N97").Value = -12390.5
$ws.Range("H131").Value = 758.3099999999999
$ws.Range("I131").Value = 266.33334
$ws.Range("J131").Value = 773.52576
$ws.Range("K131").Value = 799.0000200000001
$ws.Range("L131").Value = 2320.57728
$ws.Range("M131").Value = 4240.99998
$ws.Range("N131").Value = -12400.57728
$ws.Range("H132").Value = 1197.15
$ws.Range("I132").Value = 500
$ws.Range("J132").Value = 1233.8422
$ws.Range("K132").Value = 4500
$ws.Range("L132").Value = 11104.5798
$ws.Range("M132").Value = -1970
$ws.Range("N132").Value = -16164.5798

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 88890430
$ws.Range("I122").Value = 37037990
$ws.Range("J122").Value = 166669090
$ws.Range("K122").Value = 111113970
$ws.Range("L122").Value = 500007270
$ws.Range("M122").Value = -111111520
$ws.Range("N122").Value = -500012170
$ws.Range("H132").Value = 18232.03
$ws.Range("I132").Value = 3662.6296
$ws.Range("K132").Value = 10987.8888
$ws.Range("M132").Value = -8457.888800000001
$ws.Range("H138").Value = 50429
$ws.Range("J138").Value = 50429
$ws.Range("L138").Value = 50429
$ws.Range("N138").Value = -60709

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 2936
$ws.Range("J14").Value = 2936
$ws.Range("L14").Value = 2936
$ws.Range("N14").Value = -3280
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H132").Value = 1637.3422
$ws.Range("I132").Value = 1147.75
$ws.Range("J132").Value = 2476.6428
$ws.Range("K132").Value = 3443.25
$ws.Range("L132").Value = 7429.928400000001
$ws.Range("M132").Value = -913.25
$ws.Range("N132").Value = -12489.9284
$ws.Range("H136").Value = 21761.875
$ws.Range("I136").Value = 28582.5
$ws.Range("J136").Value = 1300
$ws.Range("K136").Value = 85747.5
$ws.Range("L136").Value = 3900
$ws.Range("M136").Value = -83197.5
$ws.Range("N136").Value = -9000

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 892.4681
$ws.Range("I132").Value = 621.2895
$ws.Range("J132").Value = 2037.4445
$ws.Range("K132").Value = 1863.8685
$ws.Range("L132").Value = 6112.333500000001
$ws.Range("M132").Value = 666.1315
$ws.Range("N132").Value = -11172.3335
$ws.Range("H136").Value = 40002092
$ws.Range("I136").Value = 55557556
$ws.Range("J136").Value = 2329.8572
$ws.Range("K136").Value = 166672668
$ws.Range("L136").Value = 6989.571599999999
$ws.Range("M136").Value = -166670118
$ws.Range("N136").Value = -12089.5716
